$d = $word.ActiveDocument
$newText = @"
# Sambar Recipe
Sambar is a popular South Indian dish, often served with rice, dosa, idli or vada. It's a comforting and flavourful lentil stew, made with vegetables, tamarind and spices. There are many different ways to prepare sambar, but here's a recipe that uses the pressure cooker method which takes around 45 minutes in total.

## Ingredients:
For the lentils:
- ¾ cup toor dal (you can also use moong dal or a combination of both)
- 2 cups water

For the vegetables:
- ½ an onion, chopped
- 1 carrot, chopped
- 5 beans, chopped
- 2 brinjal (aubergine), chopped
- 1 small eggplant, cut into 2-inch pieces
- 1 medium potato, peeled and cut into 2-inch pieces
- 1 tomato, chopped
- 1 small cucumber, chopped
- 1 green chilli, slit
- 5-6 curry leaves
- 2½ cups water

For the spice mix:
- 2 tbsp sambar powder (you can use a store-bought version or make your own at home)
- ½ tsp turmeric
- 1 tsp jaggery
- 1 tbsp tamarind extract (or the pulp of 1-2 tablespoons of tamarind soaked in hot water)
- 1 tsp salt

For the tempering:
- 2 tsp oil (ghee or vegetable oil will work)
- 1 tsp mustard
- 1 tsp urad dal
- 2 dried red chillies
- a pinch of hing (asafoetida)
- 5-6 curry leaves

For garnish:
- Fresh coriander leaves

## Method:
1. Add the lentils, a pinch of turmeric and 1 teaspoon of oil to a pressure cooker along with 2 cups of water. Cook for around 5 whistles or until the lentils are soft and mushy.
2. In a separate pot, heat some oil and add the vegetables, curry leaves and green chilli. Cook for around 10-12 minutes until the vegetables are tender but still hold their shape.
3. Add the spice powders to the pot with the vegetables and stir. 
4. Pour in the cooked lentils, along with 1½ cups of water and stir everything well.
5. Add the tamarind extract, jaggery and salt. Taste and adjust the seasoning if required. 
6. Bring the mixture to a boil and then simmer for around 3-5 minutes. 
7. Meanwhile, prepare the tempering by heating the oil in a small pan. Add the mustard seeds, urad dal, dried red chillies, hing and curry leaves. Allow them to splutter for a few seconds.
8. Pour the tempering over the sambar and mix.
9. Garnish with corriander leaves and serve hot.
"@
# Here-strings use `n` as the line separator; Word represents manual line
# breaks (<w:br/>) as vertical-tab (chr 11) in Range.Text, so convert.
$newText = $newText.Replace("`n", "`v")
$d.Content.Text = $newText
